# store.xlsx edit: add an "Other" store with a new associated value
# ("Sally's Beauty") as a third column next to the existing Store/Value
# pairs, and move the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column C holds a second value per store-row; give it a sensible
# width (closest the host lets us set explicitly) like the source file.
$ws.Columns.Item(3).ColumnWidth = 12

# Row 2 / Row 3 gain a new value in column C.
$ws.Range("C2").Value = "Other"
$ws.Range("C3").Value = "Sally's Beauty"

# Move the active selection to where the user ended up after editing.
$ws.Range("D20").Select()
